# Updating filtered feeds from workflow
# Adds two new rows for the "Foundation Medicine, Sumitomo Pharma to Develop CDx
# for Acute Leukemia Treatment" article (one link from genomeweb.com, one from
# 360dx.com), each tagged with the "CDx" keyword.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$link1 = "https://www.genomeweb.com/cancer/foundation-medicine-sumitomo-pharma-develop-cdx-acute-leukemia-treatment"
$link2 = "https://www.360dx.com/cancer/foundation-medicine-sumitomo-pharma-develop-cdx-acute-leukemia-treatment"
$keyword = "CDx"
$title = "Foundation Medicine, Sumitomo Pharma to Develop CDx for Acute Leukemia Treatment"

# Row 18
$ws.Range("A18").Value = $link1
$ws.Range("B18").Value = $keyword
$ws.Range("C18").Value = $title

# Row 19
$ws.Range("A19").Value = $link2
$ws.Range("B19").Value = $keyword
$ws.Range("C19").Value = $title

# Hyperlinks for the new link cells (matches the style of the existing rows)
$ws.Hyperlinks.Add($ws.Range("A18"), $link1)
$ws.Hyperlinks.Add($ws.Range("A19"), $link2)

# Re-apply the workbook's Hyperlink cell style so no stray style entries are used
$ws.Range("A18").Style = "Hyperlink"
$ws.Range("A19").Style = "Hyperlink"
